$wb = $excel.ActiveWorkbook

# The test data's saved-place strings are being trimmed down to the
# shorter "<street>, <city>, <state>" form (dropping the placeholder
# "5521 Alton" / "14201 Jeffrey" fragments and the full postal address
# that used to live in the sheet).

# "addPlaceTest" : 14201 Jeffrey -> 14201 Jeffrey Rd, Irvine, CA
$wsAddPlace = $wb.Worksheets.Item("addPlaceTest")
$wsAddPlace.Range("C2").Value = "14201 Jeffrey Rd, Irvine, CA"

# "addPlaceTest (OLD)" : 5521 Alton -> 5521 Alton Pkwy, Irvine, CA
$wsAddPlaceOld = $wb.Worksheets.Item("addPlaceTest (OLD)")
$wsAddPlaceOld.Activate() | Out-Null
$wsAddPlaceOld.Range("C2").Value = "5521 Alton Pkwy, Irvine, CA"
$wsAddPlaceOld.Range("C2").Select() | Out-Null

# "addLBAlertWithPlaceTest" : 14201 Jeffrey -> 14201 Jeffrey Rd, Irvine, CA
$wsAlert = $wb.Worksheets.Item("addLBAlertWithPlaceTest")
$wsAlert.Activate() | Out-Null
$wsAlert.Range("C2").Value = "14201 Jeffrey Rd, Irvine, CA"
$wsAlert.Range("C2").Select() | Out-Null

# Restore the workbook's originally active sheet/tab.
$wsAddPlace.Activate() | Out-Null
